$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 7/8: TX interval reading bumped (tics/ms recompute automatically) ---
$ws.Range("B7").Value = 34

# --- Row 9: quick packet-duration scratch calc in column F ---
$ws.Range("F9").Formula = "=23*1.5"

# --- Row 10: Tick count promoted to "Good" styling, new measured value ---
$ws.Range("A10:B10").Style = "Good"
$ws.Range("B10").Value = 5040
$ws.Range("D10").Value = "Pkt Dur"

# --- Row 11: Time promoted to "Calculation" styling (formula recalculates) ---
$ws.Range("A11:B11").Style = "Calculation"

# --- Row 13/14: second Tick-count-derived pair (Cycle duration / Time) ---
$ws.Range("A13").Value = "Cycle duration"
$ws.Range("B13").Formula = "= 1.1*B10*10"

$ws.Range("A14").Value = "Time"
$ws.Range("A14").Style = "Calculation"
$ws.Range("B14").Formula = "=B13*B5/1000"
$ws.Range("B14").Style = "Calculation"

# --- Rows 17-21: new battery/ADC measurement block ---
$ws.Range("A17").Value = "Uref"
$ws.Range("A17").Style = "Check Cell"
$ws.Range("B17").Value = 2560
$ws.Range("B17").Style = "Check Cell"
$ws.Range("C17").Value = "mV"

$ws.Range("A18").Value = "1 bit"
$ws.Range("B18").Formula = "=B17/1024"
$ws.Range("C18").Value = "mV"
$ws.Range("A18:B18").Style = "Neutral"
$ws.Range("A18").Borders.Color = 8355711
$ws.Range("A18").Borders.Weight = 2
$ws.Range("A18").Borders.LineStyle = 1
$ws.Range("B18").Borders.Color = 8355711
$ws.Range("B18").Borders.Weight = 2
$ws.Range("B18").Borders.LineStyle = 1

$ws.Range("A19").Value = "Uin"
$ws.Range("A19").Style = "Good"
$ws.Range("B19").Value = 3600
$ws.Range("B19").Style = "Good"
$ws.Range("C19").Value = "mV"

$ws.Range("A20").Value = "UinADC"
$ws.Range("A20").Style = "Neutral"
$ws.Range("B20").Formula = "=B19/2"
$ws.Range("B20").Style = "Neutral"
$ws.Range("C20").Value = "mV"

$ws.Range("A21").Value = "ADC out"
$ws.Range("A21").Style = "Calculation"
$ws.Range("B21").Formula = "=B20/B18"
$ws.Range("B21").Style = "Calculation"

# --- Selection left where the author's cursor ended up ---
$null = $ws.Range("G22").Select()
